$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit after the
#    "Was not the best father..." paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Drunk man's "Name: " paragraph gains a second run containing
#    "Dude".
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Name: `r") {
        $r = $p.Range
        $insPoint = $r.Duplicate
        $insPoint.Collapse(0)
        [void]$insPoint.MoveEnd(1, -1)
        $startPos = $insPoint.End
        $insPoint.InsertAfter("X")
        $newRange = $d.Range($startPos, $startPos + 1)
        $d.Bookmarks.Add("TMP_DUDE", $newRange)
        $bmr = $d.Bookmarks("TMP_DUDE").Range
        $bmr.Text = "Dude"
        $d.Bookmarks("TMP_DUDE").Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3) "Name: Kira (killer)" -> split "killer" into its own run and
#    wrap it with a (newly recreated) "_GoBack" bookmark.
# ------------------------------------------------------------------
$killerRange = $d.Content
$found = $killerRange.Find.Execute("killer", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if ($found) {
    $d.Bookmarks.Add("_GoBack", $killerRange)
}
